$rowsData = @(
    "142|44800|2252.4587351766299|2225.5|0|191|",
    "143|44800|2249.5919395711398|2224.4699999999998|3|180|CRM opened 8/21/2022",
    "144|44800|2229.9277457974299|2230.52|0|183|New CRM opened 8/27/2022",
    "145|44803|2232.5404066472402|2230.52|0|183|CRM opened 8/27/2022",
    "146|44805|2246.7591080000002|2230.52|0|183|CRM opened 8/27/2022",
    "147|44882|2251.5197182939401|2225.5|3|191|CRM opened 11/17/2022",
    "148|44882|2247.6284427516698|2225.5|3|191|CRM opened 11/17/2022",
    "149|44882|2236.2906877866599|2225.5|3|191|CRM opened 11/17/2022",
    "150|44882|2221.0105802766798|2225.5|3|191|CRM opened 11/17/2022",
    "151|44883|2240.1143773028998|2225.5|3|191|CRM opened 11/17/2022",
    "152|44908|2253.5423536174399|2225.5|3|191|CRM opened 11/17/2022",
    "153|44908|2248.8867845742002|2225.5|3|191|CRM opened 11/17/2022",
    "154|44908|2223.8078839847099|2225.5|3|191|CRM opened 2022-12-13",
    "155|44910|2236.8053891996001|2225.5|3|191|CRM opened 2022-12-13",
    "156|44911|2240.1060148698698|2225.5|3|191|CRM opened 2022-12-13",
    "157|44914|2224.52623796326|2225.5|3|191|CRM opened 2022-12-13",
    "158|44915|2249.06168130893|2225.5|3|191|CRM opened 2022-12-13",
    "159|44915|2232.1446972612798|2225.5|3|191|CRM opened 2022-12-13",
    "160|44916|2222.9937336656999|2225.5|3|191|CRM opened 2022-12-21",
    "161|45021|2202.2406761348202|2215.13|3|202|CRM opened 4/4/2023",
    "162|45022|2230.1596600173102|2215.13|3|202|CRM opened 4/4/2023",
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number/font formats from existing template cells so no new style
# entries are introduced (exactly mirrors existing style indices 1 and 3).
$ws.Range("A141").Copy()
$ws.Range("A142:A162").PasteSpecial(-4122)

$ws.Range("C141").Copy()
$ws.Range("C143").PasteSpecial(-4122)
$ws.Range("C147:C162").PasteSpecial(-4122)

foreach ($line in $rowsData) {
    $parts = $line.Split("|")
    $r      = [int]$parts[0]
    $aVal   = [double]$parts[1]
    $bVal   = [double]$parts[2]
    $cVal   = [double]$parts[3]
    $eVal   = [int]$parts[5]
    $fVal   = $parts[6]

    $ws.Cells.Item($r, 1).Value = $aVal
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Formula = "=100*(B$r-C$r)/C$r"
    $ws.Cells.Item($r, 5).Value = $eVal
    if ($fVal -ne "") {
        $ws.Cells.Item($r, 6).Value = $fVal
    }
}

# Restore the view state to match the edit: zoom to 130% and move the
# selection down to the new last row (freeze stays ySplit=1, already set).
$excel.ActiveWindow.Zoom = 130
$ws.Range("D164").Select()

Write-Host "done"
